$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 3906.399109145206
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 48353.76274462014
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 9433.134471502228
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2534.277928792104
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2367.37219622158
$ws.Range("O2").Value = 1995.762462679798

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 6991.052031681918
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 197913.7502057619
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16452.51445364119
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 8194.52068131253
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 7543.193583625169
$ws.Range("O2").Value = 6257.586732772244

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 12888.04225687751
$ws.Range("O2").Value = 9263.466444480218

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 31236.29455387744
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 14045.89200932069
$ws.Range("O2").Value = 9263.466444480218

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 16879.89729726143
$ws.Range("O2").Value = 10096.02314047837

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("A1").Value = "eb"
$ws.Range("B1").Value = "gb"
$ws.Range("C1").Value = "hp"
$ws.Range("D1").Value = "st"
$ws.Range("E1").Value = "wi"
$ws.Range("F1").Value = "ieh"
$ws.Range("G1").Value = "chp"
$ws.Range("H1").Value = "ac"
$ws.Range("I1").Value = "ab_ct"
$ws.Range("J1").Value = "ab_hp"
$ws.Range("K1").Value = "cp_ct"
$ws.Range("L1").Value = "cp_hp"
$ws.Range("M1").Value = "ttes"
$ws.Range("N1").Value = "btes"
$ws.Range("O1").Value = "ites"
$ws.Range("A2").Value = 38906.8534480406
$ws.Range("B2").Value = 193.0947398408091
$ws.Range("C2").Value = 292247.2772138842
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 289724.0114301849
$ws.Range("F2").Value = 16595.10705160327
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12131.91920790125
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 16879.89729726143
$ws.Range("O2").Value = 10096.02314047837
